$wb = $excel.ActiveWorkbook

# Rename Sheet1 to "Configuration Testing"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Configuration Testing"

# Update window size (zoom/resize the workbook window)
$excel.ActiveWindow.Width = 17175
$excel.ActiveWindow.Height = 11130

# Fill in new text for B2/C2 and B3/C3, with row height 30 and text wrapping already set by style
$ws.Range("B2").Value = "Go into the LbcbPlugin Folder and double click on the LbcbPlugin code file"
$ws.Range("C2").Value = "MATLAB should start up with a command window and an editor window"
$ws.Range("C3").Value = "The window should clear."
$ws.Range("B3").Value = "You can dismiss the editor window.  Type ""clearSpace"" in the command window.  "

# Set row heights for rows 2 and 3
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# Update the active selection to B3
$ws.Range("B3").Select()
